# Generate Report for Handback
# Adds a new handback entry (42de6f1f-6151-4387-8608-f2709067f14e) as row 4
# on the "Overview", "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$fileId   = "42de6f1f-6151-4387-8608-f2709067f14e"
$mdName   = "$fileId.md"
$zhXlf    = "$fileId.ad411bf5d8a339d31bdea6f1ed01c9b543ef7ef2.zh-cn.xlf"
$deXlf    = "$fileId.ad411bf5d8a339d31bdea6f1ed01c9b543ef7ef2.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"
$includeText  = "Include"

# -----------------------------------------------------------------
# Sheet 1: "Overview"
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Cells.Item(4, 2).Value = $statusInSync
$ws1.Cells.Item(4, 3).Value = $statusInSync

$ws1.Hyperlinks.Add(
    $ws1.Cells.Item(4, 1),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/82a3b8f1e3940651f21c6f0e2393c586908da162/e2e/$mdName",
    $null,
    $null,
    $mdName
)

# -----------------------------------------------------------------
# Sheet 2: "zh-cn"
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Cells.Item(4, 2).Value = $statusInSync
$ws2.Cells.Item(4, 4).Value = "2016-03-04 08:19:03"
$ws2.Cells.Item(4, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(4, 7).Value = "2016-03-04 08:19:46"
$ws2.Cells.Item(4, 8).Value = $includeText

$ws2.Hyperlinks.Add(
    $ws2.Cells.Item(4, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/82a3b8f1e3940651f21c6f0e2393c586908da162/e2e/$mdName",
    $null,
    $null,
    $mdName
)
$ws2.Hyperlinks.Add(
    $ws2.Cells.Item(4, 3),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f5c5af615442aa7f666331bd3d956ca789976ef0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf",
    $null,
    $null,
    $zhXlf
)
$ws2.Hyperlinks.Add(
    $ws2.Cells.Item(4, 5),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a808c91c1cb7f6b5df575e14cf0ba8ba48a024ff/e2e/$mdName",
    $null,
    $null,
    $mdName
)
$ws2.Hyperlinks.Add(
    $ws2.Cells.Item(4, 6),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/be87076052a3fdc54956cd2cdbcfdd3e19bd0cf7/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf",
    $null,
    $null,
    $zhXlf
)

# -----------------------------------------------------------------
# Sheet 3: "de-de"
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Cells.Item(4, 2).Value = $statusInSync
$ws3.Cells.Item(4, 4).Value = "2016-03-04 08:19:13"
$ws3.Cells.Item(4, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(4, 7).Value = "2016-03-04 08:20:05"
$ws3.Cells.Item(4, 8).Value = $includeText

$ws3.Hyperlinks.Add(
    $ws3.Cells.Item(4, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/82a3b8f1e3940651f21c6f0e2393c586908da162/e2e/$mdName",
    $null,
    $null,
    $mdName
)
$ws3.Hyperlinks.Add(
    $ws3.Cells.Item(4, 3),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/900391d25e450fb0527033c5942a70618d2c0b40/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf",
    $null,
    $null,
    $deXlf
)
$ws3.Hyperlinks.Add(
    $ws3.Cells.Item(4, 5),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0e6b2395c20fbee4147e141be1114f02d10a0a00/e2e/$mdName",
    $null,
    $null,
    $mdName
)
$ws3.Hyperlinks.Add(
    $ws3.Cells.Item(4, 6),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/651bc3af75781e154d52c74adbb57fcfe8780a63/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf",
    $null,
    $null,
    $deXlf
)
